$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.370.26'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.874.97'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7119'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.05'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07815'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3121'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.19'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08434'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.875.86'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7136'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.07'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.380.51'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.056'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008237'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +5.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.81'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.25'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.122.19'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.787'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1592'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.99%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.57'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.511'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.290'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.92%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.324'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05299'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.939'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.179'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7445'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -9.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.700'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01870'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.227.27'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.67%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.29%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +4.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '110.81'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +8.56%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.80'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.020.21'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.812'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5212'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000123'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.428'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4327'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.36%  '
